# Commit: "Remove unused code and update template"
#
# tbl_spesifikasi currently has a single "os" column between "processor" and
# "memory". Replace it with three columns: "os1", "os2", "os3" (inserting two
# extra columns), and make tbl_spesifikasi the active sheet/tab instead of
# tbl_aset.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tbl_spesifikasi")

# Insert two extra columns right after the existing "os" column (column E)
# so there is room for os1/os2/os3 in place of the single os column.
$ws.Columns.Item(6).Insert()
$ws.Columns.Item(6).Insert()

# Rename the header cells: E2 was "os"; now E2/F2/G2 become os1/os2/os3.
$ws.Cells.Item(2, 5).Value = "os1"
$ws.Cells.Item(2, 6).Value = "os2"
$ws.Cells.Item(2, 7).Value = "os3"

# Make tbl_spesifikasi the active sheet/tab (was tbl_aset), with F5 selected.
[void]$ws.Activate()
[void]$ws.Range("F5").Select()
